$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2634.625
$ws.Range("I2").Value = 871.5714
$ws.Range("K2").Value = 871.5714
$ws.Range("M2").Value = -758.5714

$ws.Range("H9").Value = 126.666664
$ws.Range("I9").Value = 80
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 80
$ws.Range("L9").Value = 150
$ws.Range("M9").Value = 89
$ws.Range("N9").Value = -488

$ws.Range("H38").Value = 816.1429000000001
$ws.Range("I38").Value = 118.833336
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 356.500008
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = 15.49999200000002
$ws.Range("N38").Value = -15744

$ws.Range("H100").Value = 2886.5
$ws.Range("I100").Value = 2865.3333
$ws.Range("K100").Value = 2865.3333
$ws.Range("M100").Value = -2324.3333

$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = $null
$ws.Range("N107").Value = -80

$ws.Range("H111").Value = 487.25
$ws.Range("I111").Value = 449.66666
$ws.Range("K111").Value = 1348.99998
$ws.Range("M111").Value = 1718.00002

$ws.Range("H113").Value = 8794.786
$ws.Range("I113").Value = 6145.8
$ws.Range("J113").Value = 10266.444
$ws.Range("K113").Value = 6145.8
$ws.Range("L113").Value = 10266.444
$ws.Range("M113").Value = -2891.8
$ws.Range("N113").Value = -16774.444

$ws.Range("H116").Value = 10000
$ws.Range("I116").Value = 10000
$ws.Range("K116").Value = 10000
$ws.Range("M116").Value = -6558

$ws.Range("H118").Value = 202
$ws.Range("I118").Value = 202
$ws.Range("K118").Value = 606
$ws.Range("M118").Value = 1051

$ws.Range("H125").Value = 3664.5
$ws.Range("I125").Value = 2156.2222
$ws.Range("J125").Value = 6379.4
$ws.Range("K125").Value = 19405.9998
$ws.Range("L125").Value = 57414.6
$ws.Range("M125").Value = -16945.9998
$ws.Range("N125").Value = -62334.6

$ws.Range("H127").Value = 1800
$ws.Range("I127").Value = 1400
$ws.Range("K127").Value = 4200
$ws.Range("M127").Value = 760

$ws.Range("H132").Value = 1919.0834
$ws.Range("I132").Value = 1919.0834
$ws.Range("K132").Value = 5757.2502
$ws.Range("M132").Value = -3227.2502

$ws.Range("H137").Value = 1882.76
$ws.Range("I137").Value = 1602.55
$ws.Range("J137").Value = 3003.6
$ws.Range("K137").Value = 4807.65
$ws.Range("L137").Value = 9010.799999999999
$ws.Range("M137").Value = -2257.65
$ws.Range("N137").Value = -14110.8

$ws.Range("H138").Value = 2446.5652
$ws.Range("I138").Value = 1396.4
$ws.Range("J138").Value = 2738.2778
$ws.Range("K138").Value = 4189.200000000001
$ws.Range("L138").Value = 8214.8334
$ws.Range("M138").Value = 950.7999999999993
$ws.Range("N138").Value = -18494.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1098.8182
$ws.Range("I2").Value = 1173.625
$ws.Range("J2").Value = 899.3333
$ws.Range("K2").Value = 1173.625
$ws.Range("L2").Value = 899.3333
$ws.Range("M2").Value = -1060.625
$ws.Range("N2").Value = -1125.3333

$ws.Range("H8").Value = 10004915
$ws.Range("I8").Value = 25002284
$ws.Range("J8").Value = 6668.3335
$ws.Range("K8").Value = 25002284
$ws.Range("L8").Value = 6668.3335
$ws.Range("M8").Value = -25002140
$ws.Range("N8").Value = -6956.3335

$ws.Range("H32").Value = 1656.2142
$ws.Range("I32").Value = 1535.2898
$ws.Range("K32").Value = 1535.2898
$ws.Range("M32").Value = -1248.2898

$ws.Range("H55").Value = 27220.555
$ws.Range("J55").Value = 27220.555
$ws.Range("L55").Value = 27220.555
$ws.Range("N55").Value = -27850.555

$ws.Range("H61").Value = 8000
$ws.Range("I61").Value = 8000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 8000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("N61").Value = -7788

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = $null

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = $null

$ws.Range("H97").Value = 522.3077
$ws.Range("I97").Value = 491.25
$ws.Range("K97").Value = 491.25
$ws.Range("M97").Value = 4.75

$ws.Range("H110").Value = 873.75
$ws.Range("I110").Value = 500
$ws.Range("J110").Value = 998.3333
$ws.Range("K110").Value = 500
$ws.Range("L110").Value = 998.3333
$ws.Range("M110").Value = 1545
$ws.Range("N110").Value = -5088.3333

$ws.Range("H116").Value = 1098.8182
$ws.Range("I116").Value = 1173.625
$ws.Range("J116").Value = 899.3333
$ws.Range("K116").Value = 1173.625
$ws.Range("L116").Value = 899.3333
$ws.Range("M116").Value = 1120.375
$ws.Range("N116").Value = -5487.3333

$ws.Range("H124").Value = 40429
$ws.Range("J124").Value = 40429
$ws.Range("L124").Value = 40429
$ws.Range("N124").Value = -50249

$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 8000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 24000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = -21450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1098.8182
$ws.Range("I3").Value = 1173.625
$ws.Range("J3").Value = 899.3333
$ws.Range("K3").Value = 1173.625
$ws.Range("L3").Value = 899.3333
$ws.Range("M3").Value = -1059.625
$ws.Range("N3").Value = -1127.3333

$ws.Range("H94").Value = 1297.9166
$ws.Range("I94").Value = 1558.125
$ws.Range("J94").Value = 777.5
$ws.Range("K94").Value = 1558.125
$ws.Range("L94").Value = 777.5
$ws.Range("M94").Value = -1107.125
$ws.Range("N94").Value = -1679.5

$ws.Range("H99").Value = 2056.75
$ws.Range("I99").Value = 1779.1428
$ws.Range("K99").Value = 1779.1428
$ws.Range("M99").Value = -281.1428000000001

$ws.Range("H107").Value = 1883.1666
$ws.Range("I107").Value = 1699.75
$ws.Range("J107").Value = 2250
$ws.Range("K107").Value = 1699.75
$ws.Range("L107").Value = 2250
$ws.Range("M107").Value = 220.25
$ws.Range("N107").Value = -6090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 34857
$ws.Range("I59").Value = 33999
$ws.Range("J59").Value = 35000
$ws.Range("K59").Value = 33999
$ws.Range("L59").Value = 35000
$ws.Range("M59").Value = -32854
$ws.Range("N59").Value = -37290

$ws.Range("H86").Value = 7138.6665
$ws.Range("I86").Value = 6736.222
$ws.Range("J86").Value = 7742.3335
$ws.Range("K86").Value = 6736.222
$ws.Range("L86").Value = 7742.3335
$ws.Range("M86").Value = -5613.222
$ws.Range("N86").Value = -9988.333500000001

$ws.Range("H89").Value = 7138.6665
$ws.Range("I89").Value = 6736.222
$ws.Range("J89").Value = 7742.3335
$ws.Range("K89").Value = 33681.11
$ws.Range("L89").Value = 38711.6675
$ws.Range("M89").Value = -28065.11
$ws.Range("N89").Value = -49943.6675

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 900
$ws.Range("I131").Value = 800
$ws.Range("K131").Value = 2400
$ws.Range("M131").Value = 2640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 6373.6
$ws.Range("I99").Value = 6373.6
$ws.Range("K99").Value = 6373.6
$ws.Range("M99").Value = -4127.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17773.8
$ws.Range("I7").Value = 18472.422
$ws.Range("K7").Value = 18472.422
$ws.Range("M7").Value = -18360.422

$ws.Range("H26").Value = 50010
$ws.Range("J26").Value = 50010
$ws.Range("L26").Value = 50010
$ws.Range("N26").Value = -50600

$ws.Range("H53").Value = 32500
$ws.Range("I53").Value = 32500
$ws.Range("K53").Value = 32500
$ws.Range("M53").Value = -31982

$ws.Range("H93").Value = 1424.5
$ws.Range("I93").Value = 1360.5555
$ws.Range("K93").Value = 1360.5555
$ws.Range("M93").Value = -112.5554999999999

$ws.Range("H95").Value = 29124.5
$ws.Range("J95").Value = 29124.5
$ws.Range("L95").Value = 29124.5
$ws.Range("N95").Value = -34616.5

$ws.Range("H100").Value = 4499.875
$ws.Range("J100").Value = 4999.857
$ws.Range("L100").Value = 4999.857
$ws.Range("N100").Value = -6081.857

$ws.Range("H104").Value = 38074
$ws.Range("J104").Value = 38074
$ws.Range("L104").Value = 38074
$ws.Range("N104").Value = -45062

$ws.Range("H126").Value = 17773.8
$ws.Range("I126").Value = 18472.422
$ws.Range("K126").Value = 55417.266
$ws.Range("M126").Value = -52947.266

$ws.Range("H136").Value = 2817.8667
$ws.Range("I136").Value = 2635.6155
$ws.Range("K136").Value = 7906.8465
$ws.Range("M136").Value = -5356.8465

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 1000
$ws.Range("L10").Value = 1000
$ws.Range("N10").Value = -1338

$ws.Range("H13").Value = 7863.75
$ws.Range("J13").Value = 10151.667
$ws.Range("L13").Value = 10151.667
$ws.Range("N13").Value = -10431.667

$ws.Range("H97").Value = 6250
$ws.Range("J97").Value = 6250
$ws.Range("L97").Value = 6250
$ws.Range("N97").Value = -8232
